$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1).
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 22:14"

# Full refreshed country table (country name + 7 numeric stat columns),
# already sorted by "Casos totales" (column B) descending, matching the
# post-refresh row order (several countries changed rank).
$data = @(
  ,@("Estados Unidos",3668280,51453,1663668,1863827,0,641,140785)
  ,@("Brasil",1978236,7327,1366775,535764,0,174,75697)
  ,@("India",1005637,35468,636602,343426,0,680,25609)
  ,@("Rusia",752797,6428,531692,209168,0,167,11937)
  ,@("Peru",337724,0,226400,98907,0,0,12417)
  ,@("Sudafrica",324221,13172,165591,153961,0,216,4669)
  ,@("Chile",323698,2493,295301,21107,0,104,7290)
  ,@("Mexico",317635,6149,199129,81600,0,579,36906)
  ,@("España",305935,1361,0,0,0,3,28416)
  ,@("Reino Unido",292552,641,0,0,0,66,45119)
  ,@("Iran",267061,2500,230608,22845,0,198,13608)
  ,@("Pakistan",257914,2145,178737,73751,0,40,5426)
  ,@("Italia",243736,230,196246,12473,0,20,35017)
  ,@("Arabia Saudita",243238,2764,187622,53246,0,45,2370)
  ,@("Turquia",216873,933,198820,12613,0,21,5440)
  ,@("Alemania",201827,575,186400,6270,0,9,9157)
  ,@("Banglades",196323,2733,106963,86864,0,39,2496)
  ,@("Francia",173838,534,79036,64664,0,18,30138)
  ,@("Colombia",165169,0,71736,87619,0,0,5814)
  ,@("Argentina",111146,0,47298,61798,0,0,2050)
  ,@("Canada",109080,251,72782,27473,0,15,8825)
  ,@("Catar",105477,494,102168,3157,0,1,152)
  ,@("Irak",86148,2281,54316,28310,0,90,3522)
  ,@("Egipto",84843,0,26135,54641,0,0,4067)
  ,@("China",83612,1,78719,259,0,0,4634)
  ,@("Indonesia",81668,1574,40345,37450,0,76,3873)
  ,@("Suecia",76877,115,0,0,0,10,5593)
  ,@("Ecuador",71365,1036,31260,34898,0,49,5207)
  ,@("Bielorrusia",65623,180,56915,8223,0,5,485)
  ,@("Kazajistan",65188,1674,39066,25747,0,0,375)
  ,@("Belgica",63039,167,17253,35994,0,4,9792)
  ,@("Oman",62574,1327,40090,22194,0,9,290)
  ,@("Filipinas",61266,2498,21440,38183,0,29,1643)
  ,@("Kuwait",57668,791,47545,9721,0,3,402)
  ,@("Ucrania",56455,848,28931,26079,0,18,1445)
  ,@("Emiratos Arabes Unidos",56129,281,47412,8382,0,0,335)
  ,@("Bolivia",52218,1351,16357,33919,0,44,1942)
  ,@("Paises Bajos",51351,99,0,0,0,1,6137)
  ,@("Panama",49243,0,25417,22844,0,0,982)
  ,@("Republica Dominicana",48743,1072,23636,24166,0,12,941)
  ,@("Portugal",47765,339,32476,13610,0,3,1679)
  ,@("Singapur",47126,248,43256,3843,0,0,27)
  ,@("Israel",45607,1419,20268,24956,0,7,383)
  ,@("Polonia",39054,333,28928,8521,0,11,1605)
  ,@("Afganistan",35070,76,22824,11133,0,19,1113)
  ,@("Rumania",35003,777,22189,10843,0,19,1971)
  ,@("Barein",34560,0,30320,4119,0,4,121)
  ,@("Nigeria",34259,0,13999,19500,0,0,760)
  ,@("Armenia",33559,554,21931,11021,0,15,607)
  ,@("Suiza",33290,142,29800,1521,0,1,1969)
  ,@("Guatemala",32074,0,4624,26100,0,0,1350)
  ,@("Honduras",30036,930,3379,25832,0,18,825)
  ,@("Azerbaiyan",26165,493,17256,8575,0,8,334)
  ,@("Ghana",26125,695,22270,3716,0,0,139)
  ,@("Irlanda",25698,15,23364,585,0,1,1749)
  ,@("Japon",22508,0,18545,2979,0,0,984)
  ,@("Argelia",21355,585,15107,5196,0,12,1052)
  ,@("Moldavia",20264,224,13640,5958,0,7,666)
  ,@("Serbia",19717,383,14047,5228,0,13,442)
  ,@("Austria",19270,116,17244,1315,0,1,711)
  ,@("Nepal",17344,167,11249,6056,0,0,39)
  ,@("Marruecos",16545,283,13965,2317,0,4,263)
  ,@("Camerun",15173,0,11928,2886,0,0,359)
  ,@("Uzbekistan",14991,410,8655,6263,0,2,73)
  ,@("Corea del Sur",13612,61,12396,925,0,2,291)
  ,@("Chequia",13551,76,8507,4689,0,0,355)
  ,@("Costa de Marfil",13403,0,7146,6170,0,0,87)
  ,@("Dinamarca",13124,32,12195,319,0,0,610)
  ,@("Kirguistan",12498,521,3735,8596,0,7,167)
  ,@("Kenia",11673,421,3638,7818,0,8,217)
  ,@("El Salvador",10957,312,6257,4402,0,12,298)
  ,@("Australia",10810,323,8035,2662,0,2,113)
  ,@("Sudan",10527,0,5601,4258,0,0,668)
  ,@("Venezuela",10428,0,3050,7278,0,0,100)
  ,@("Costa Rica",9546,560,2673,6831,0,2,42)
  ,@("Noruega",9011,0,8138,619,0,1,254)
  ,@("Malasia",8737,3,8538,77,0,0,122)
  ,@("Republica de Macedonia",8623,93,4565,3657,0,8,401)
  ,@("Senegal",8481,112,5735,2590,0,3,156)
  ,@("Consejo Danes para los Refugiados",8199,36,4248,3758,0,1,193)
  ,@("Etiopia",8181,0,2430,5605,0,0,146)
  ,@("Bulgaria",7877,0,3841,3747,0,0,289)
  ,@("Bosnia y Herzegovina",7681,270,3534,3907,0,5,240)
  ,@("Estado de Palestina",7412,348,1313,6048,0,7,51)
  ,@("Finlandia",7293,0,6880,85,0,0,328)
  ,@("Haiti",6902,71,3484,3273,0,2,145)
  ,@("Tayikistan",6695,0,5383,1256,0,0,56)
  ,@("Guayana Francesa",6393,94,3818,2541,0,1,34)
  ,@("Guinea",6276,0,4981,1257,0,0,38)
  ,@("Gabon",6121,0,3664,2411,0,0,46)
  ,@("Madagascar",6089,484,2951,3085,0,10,53)
  ,@("Mauritania",5564,0,2830,2585,0,0,149)
  ,@("Luxemburgo",5122,0,4247,764,0,0,111)
  ,@("Republica de Yibuti",4985,0,4765,164,0,0,56)
  ,@("Republica de Africa Central",4362,0,1261,3048,0,0,53)
  ,@("Hungria",4279,16,3156,528,0,0,595)
  ,@("Croacia",4039,86,2729,1190,0,0,120)
  ,@("Grecia",3939,29,1374,2372,0,0,193)
  ,@("Albania",3851,99,2137,1610,0,3,104)
  ,@("Tailandia",3236,4,3095,83,0,0,58)
  ,@("Paraguay",3198,0,1338,1835,0,0,25)
  ,@("Nicaragua",3147,0,2282,766,0,0,99)
  ,@("Somalia",3083,0,1425,1565,0,0,93)
  ,@("Guinea Ecuatorial",3071,0,842,2178,0,0,51)
  ,@("Maldivas",2831,0,2321,495,0,1,15)
  ,@("Mayotte",2743,0,2581,125,0,0,37)
  ,@("Sri Lanka",2686,15,2007,668,0,0,11)
  ,@("Malaui",2614,0,1005,1566,0,0,43)
  ,@("Libano",2599,57,1485,1074,0,2,40)
  ,@("Mali",2440,7,1777,542,0,0,121)
  ,@("Cuba",2440,2,2285,68,0,0,87)
  ,@("Congo",2222,0,589,1586,0,0,47)
  ,@("Sudan del Sur",2153,0,1175,937,0,0,41)
  ,@("Estonia",2016,0,1904,43,0,0,69)
  ,@("Eslovaquia",1951,24,1514,409,0,0,28)
  ,@("Islandia",1914,3,1892,12,0,0,10)
  ,@("Guinea-Bisau",1902,60,773,1103,0,0,26)
  ,@("Lituania",1902,6,1593,230,0,0,79)
  ,@("Eslovenia",1897,19,1522,264,0,0,111)
  ,@("Zambia",1895,0,1412,441,0,0,42)
  ,@("Cabo Verde",1837,57,892,926,0,0,19)
  ,@("Sierra Leona",1668,0,1200,404,0,0,64)
  ,@("Hong Kong",1656,67,1254,392,0,0,10)
  ,@("Libia",1652,63,379,1227,0,3,46)
  ,@("Nueva Zelanda",1548,1,1499,27,0,0,22)
  ,@("Yemen",1526,0,694,399,0,0,433)
  ,@("Suazilandia",1489,0,724,745,0,0,20)
  ,@("Ruanda",1435,0,752,679,0,0,4)
  ,@("Mozambique",1383,53,375,999,0,0,9)
  ,@("Benin",1378,0,557,795,0,0,26)
  ,@("Tunez",1327,8,1093,184,0,0,50)
  ,@("Montenegro",1287,0,330,933,0,0,24)
  ,@("Jordania",1206,5,1019,177,0,0,10)
  ,@("Letonia",1179,1,1022,126,0,0,31)
  ,@("Niger",1100,0,993,38,0,0,69)
  ,@("Zimbabue",1089,0,395,674,0,0,20)
  ,@("Liberia",1056,0,447,558,0,0,51)
  ,@("Uganda",1051,8,1014,37,0,0,0)
  ,@("Burkina Faso",1038,0,882,103,0,0,53)
  ,@("Namibia",1032,72,31,999,0,0,2)
  ,@("Republica de Chipre",1025,0,839,167,0,0,19)
  ,@("Uruguay",1009,0,909,69,0,0,31)
  ,@("Georgia",1006,2,883,108,0,0,15)
  ,@("Republica del Chad",885,0,799,11,0,0,75)
  ,@("Principado de Andorra",862,0,803,7,0,0,52)
  ,@("Surinam",849,12,573,258,0,0,18)
  ,@("Jamaica",763,1,647,106,0,0,10)
  ,@("Togo",740,0,534,191,0,0,15)
  ,@("Santo Tome y Principe",737,0,322,401,0,0,14)
  ,@("Crucero",712,0,651,48,0,0,13)
  ,@("San Marino",699,0,656,1,0,0,42)
  ,@("Malta",674,0,661,4,0,0,9)
  ,@("Reunion",612,4,472,137,0,0,3)
  ,@("Angola",576,0,124,425,0,0,27)
  ,@("Tanzania",509,0,183,305,0,0,21)
  ,@("Siria",477,19,140,315,0,0,22)
  ,@("Taiwan",452,1,440,5,0,0,7)
  ,@("Botsuana",399,0,38,360,0,0,1)
  ,@("Vietnam",381,0,356,25,0,0,0)
  ,@("Mauricio",343,0,331,2,0,0,10)
  ,@("Birmania",337,0,266,65,0,0,6)
  ,@("Isla de Man",336,0,312,0,0,0,24)
  ,@("Comoras",321,0,302,12,0,0,7)
  ,@("Guyana",313,0,156,139,0,0,18)
  ,@("Burundi",269,0,207,61,0,0,1)
  ,@("Mongolia",261,0,209,52,0,0,0)
  ,@("Lesoto",256,0,48,205,0,0,3)
  ,@("Martinica",255,0,98,142,0,0,15)
  ,@("Eritrea",232,0,107,125,0,0,0)
  ,@("Islas Caimanes",203,0,200,2,0,0,1)
  ,@("Guadalupe",190,0,157,19,0,0,14)
  ,@("Islas Feroe",188,0,188,0,0,0,0)
  ,@("Gibraltar",180,0,180,0,0,0,0)
  ,@("Camboya",166,1,133,33,0,0,0)
  ,@("Bermudas",150,0,137,4,0,0,9)
  ,@("Brunei",141,0,138,0,0,0,3)
  ,@("Trinidad yTobago",133,0,124,1,0,0,8)
  ,@("Bahamas",119,0,91,17,0,0,11)
  ,@("Monaco",109,0,98,7,0,0,4)
  ,@("Aruba",106,0,99,4,0,0,3)
  ,@("Barbados",104,0,90,7,0,0,7)
  ,@("Seychelles",100,0,27,73,0,0,0)
  ,@("Butan",84,0,78,6,0,0,0)
  ,@("Liechtenstein",84,0,81,2,0,0,1)
  ,@("Gambia",78,14,34,41,0,0,3)
  ,@("San Martin (Parte Holandesa)",78,0,63,0,0,0,15)
  ,@("Antigua y Barbuda",74,0,57,14,0,0,3)
  ,@("Islas Turcas y Caicos",72,0,12,58,0,0,2)
  ,@("Polinesia Francesa",62,0,60,2,0,0,0)
  ,@("Macao",46,0,45,1,0,0,0)
  ,@("San Martin (Parte Francesa)",43,0,37,3,0,0,3)
  ,@("Puerto Rico",39,0,1,36,0,0,2)
  ,@("Belice",39,0,21,16,0,0,2)
  ,@("San Vicente y las Granadinas",35,0,29,6,0,0,0)
  ,@("Guam",32,0,0,31,0,0,1)
  ,@("Fiyi",26,0,18,8,0,0,0)
  ,@("Curazao",26,0,24,1,0,0,1)
  ,@("Timor Oriental",24,0,24,0,0,0,0)
  ,@("Granada",23,0,23,0,0,0,0)
  ,@("Santa Lucia",22,0,19,3,0,0,0)
  ,@("Nueva Caledonia",22,0,21,1,0,0,0)
  ,@("Laos",19,0,19,0,0,0,0)
  ,@("Dominica",18,0,18,0,0,0,0)
  ,@("Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0)
  ,@("San Cristobal y Nieves",17,0,15,2,0,0,0)
  ,@("Islas Malvinas",13,0,13,0,0,0,0)
  ,@("Groenlandia",13,0,13,0,0,0,0)
  ,@("Montserrat",12,0,10,1,0,0,1)
  ,@("Santa Sede",12,0,12,0,0,0,0)
  ,@("Papua Nueva Guinea",11,0,8,3,0,0,0)
  ,@("Sahara Occidental",10,0,8,1,0,0,1)
  ,@("Bonaire, San Eustaquio y Saba",9,0,7,2,0,0,0)
  ,@("Islas Virgenes Britanicas",8,0,7,0,0,0,1)
  ,@("San Bartolome",6,0,6,0,0,0,0)
  ,@("Anguila",3,0,3,0,0,0,0)
  ,@("San Pedro y Miquelon",2,0,1,1,0,0,0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i
  $rec = $data[$i]
  $ws.Cells.Item($row, 1).Value = $rec[0]
  $ws.Cells.Item($row, 2).Value = [double]$rec[1]
  $ws.Cells.Item($row, 3).Value = [double]$rec[2]
  $ws.Cells.Item($row, 4).Value = [double]$rec[3]
  $ws.Cells.Item($row, 5).Value = [double]$rec[4]
  $ws.Cells.Item($row, 6).Value = [double]$rec[5]
  $ws.Cells.Item($row, 7).Value = [double]$rec[6]
  $ws.Cells.Item($row, 8).Value = [double]$rec[7]
}

Write-Output "Updated $($data.Count) country rows."
